# Agregada una nueva recepcionista y los datos de semana 18/02/2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: Constanza - new receptionist
$ws.Range("A6").Value = "Constanza"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# New column I: week 18_02_2024
$ws.Range("I1").Value = "18_02_2024"
$ws.Range("I2").Value = 3215
$ws.Range("I3").Value = 2640
$ws.Range("I4").Value = 3769
$ws.Range("I5").Value = 7034
$ws.Range("I6").Value = 33

# Underline the I3 cell (highlighted value) per style diff
$ws.Range("I3").Font.Underline = $true

# Page setup (paper size / orientation) as reflected by the printing metadata
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection after edits
$ws.Range("I3").Select()
